# Översikt BORGHOLM.xlsx - automatic update of files
#
# 1) Refresh the "Förändrad" (C) column for every data row (2..130) from
#    2023-09-13 (45182) to 2023-09-15 (45184).
# 2) A new logging notice "A 32640-2022" gained an extra observed species
#    ("Svart trolldruva"), which moves its record to the top of its
#    date-tied block (row 20) and pushes the records that used to occupy
#    rows 20-23 ("A 16171-2023", "A 57207-2018", "A 31733-2022",
#    "A 32174-2022") down by one row (21-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Blanket date refresh for every data row -----------------------
$ws.Range("C2:C130").Value = 45184

# --- 2) Rows 20-24: reorder + update the "A 32640-2022" record --------

# New row 20 : "A 32640-2022" (previously row 24), with the newly found
# species "Svart trolldruva" added (Signalarter 2->3, Alla arter 2->3).
$ws.Range("A20").Value = "A 32640-2022"
$ws.Range("B20").Value = 44783
$ws.Range("C20").Value = 45184
$ws.Range("D20").Value = "KALMAR LÄN"
$ws.Range("E20").Value = "BORGHOLM"
$ws.Range("F20").Value = "Sveaskog"
$ws.Range("G20").Value = 1.2
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 3
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = "Kornknutmossa`r`nMurgröna`r`nSvart trolldruva"
$ws.Range("S20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/artfynd/A 32640-2022.xlsx")'
$ws.Range("T20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/kartor/A 32640-2022.png")'
$ws.Range("V20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomål/A 32640-2022.docx")'
$ws.Range("W20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomålsmail/A 32640-2022.docx")'
$ws.Range("X20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsyn/A 32640-2022.docx")'
$ws.Range("Y20").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsynsmail/A 32640-2022.docx")'

# Row 21 : "A 16171-2023" (previously row 20) - content unchanged, just shifted.
$ws.Range("A21").Value = "A 16171-2023"
$ws.Range("B21").Value = 45027
$ws.Range("C21").Value = 45184
$ws.Range("D21").Value = "KALMAR LÄN"
$ws.Range("E21").Value = "BORGHOLM"
$ws.Range("F21").ClearContents()
$ws.Range("G21").Value = 12.2
$ws.Range("H21").Value = 1
$ws.Range("I21").Value = 2
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 3
$ws.Range("R21").Value = "Murgröna`r`nSårläka`r`nBlåsippa"
$ws.Range("S21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/artfynd/A 16171-2023.xlsx")'
$ws.Range("T21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/kartor/A 16171-2023.png")'
$ws.Range("V21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomål/A 16171-2023.docx")'
$ws.Range("W21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomålsmail/A 16171-2023.docx")'
$ws.Range("X21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsyn/A 16171-2023.docx")'
$ws.Range("Y21").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsynsmail/A 16171-2023.docx")'

# Row 22 : "A 57207-2018" (previously row 21) - content unchanged, just shifted.
$ws.Range("A22").Value = "A 57207-2018"
$ws.Range("B22").Value = 43403
$ws.Range("C22").Value = 45184
$ws.Range("D22").Value = "KALMAR LÄN"
$ws.Range("E22").Value = "BORGHOLM"
$ws.Range("F22").ClearContents()
$ws.Range("G22").Value = 3.6
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 1
$ws.Range("P22").Value = 1
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = "Luddvicker`r`nSkogsknipprot"
$ws.Range("S22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/artfynd/A 57207-2018.xlsx")'
$ws.Range("T22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/kartor/A 57207-2018.png")'
$ws.Range("V22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomål/A 57207-2018.docx")'
$ws.Range("W22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomålsmail/A 57207-2018.docx")'
$ws.Range("X22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsyn/A 57207-2018.docx")'
$ws.Range("Y22").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsynsmail/A 57207-2018.docx")'

# Row 23 : "A 31733-2022" (previously row 22) - content unchanged, just shifted.
$ws.Range("A23").Value = "A 31733-2022"
$ws.Range("B23").Value = 44776
$ws.Range("C23").Value = 45184
$ws.Range("D23").Value = "KALMAR LÄN"
$ws.Range("E23").Value = "BORGHOLM"
$ws.Range("F23").Value = "Sveaskog"
$ws.Range("G23").Value = 6.9
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 1
$ws.Range("P23").Value = 1
$ws.Range("Q23").Value = 2
$ws.Range("R23").Value = "Ask`r`nSkogsknipprot"
$ws.Range("S23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/artfynd/A 31733-2022.xlsx")'
$ws.Range("T23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/kartor/A 31733-2022.png")'
$ws.Range("V23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomål/A 31733-2022.docx")'
$ws.Range("W23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomålsmail/A 31733-2022.docx")'
$ws.Range("X23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsyn/A 31733-2022.docx")'
$ws.Range("Y23").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsynsmail/A 31733-2022.docx")'

# Row 24 : "A 32174-2022" (previously row 23) - content unchanged, just shifted.
$ws.Range("A24").Value = "A 32174-2022"
$ws.Range("B24").Value = 44781
$ws.Range("C24").Value = 45184
$ws.Range("D24").Value = "KALMAR LÄN"
$ws.Range("E24").Value = "BORGHOLM"
$ws.Range("F24").Value = "Sveaskog"
$ws.Range("G24").Value = 1.6
$ws.Range("H24").Value = 1
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = 1
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 1
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 2
$ws.Range("R24").Value = "Backklöver`r`nTvåblad"
$ws.Range("S24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/artfynd/A 32174-2022.xlsx")'
$ws.Range("T24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/kartor/A 32174-2022.png")'
$ws.Range("V24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomål/A 32174-2022.docx")'
$ws.Range("W24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/klagomålsmail/A 32174-2022.docx")'
$ws.Range("X24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsyn/A 32174-2022.docx")'
$ws.Range("Y24").Formula = '=HYPERLINK("https://klasma.github.io/Logging_BORGHOLM/tillsynsmail/A 32174-2022.docx")'

# Re-entering the wrapped "Artnamn" text above made the engine auto-fit
# these rows' height; restore the original fixed row height (15).
$ws.Rows("20:24").RowHeight = 15
